$wb = $excel.ActiveWorkbook

# The "max-arrecad" sheet (column A: categoria_mencao) needs to be
# re-fixed so that categories tied on the same arrecadado value keep a
# stable/deterministic order (fixing the author+mencao DB to the 2024
# analysis year).
$wsMax = $wb.Worksheets.Item("max-arrecad")

$wsMax.Range("A2").Value  = "humor"
$wsMax.Range("A3").Value  = "folclore"
$wsMax.Range("A4").Value  = "religiosidade"
$wsMax.Range("A5").Value  = "terror"
$wsMax.Range("A6").Value  = "fiq"
$wsMax.Range("A7").Value  = "ficcao_cientifica"
$wsMax.Range("A8").Value  = "fantasia"
$wsMax.Range("A9").Value  = "questoes_genero"

$wsMax.Range("A11").Value = "jogos"
$wsMax.Range("A12").Value = "webformatos"

$wsMax.Range("A13").Value = "hqmix"
$wsMax.Range("A14").Value = "angelo_agostini"

$wsMax.Range("A16").Value = "erotismo"
$wsMax.Range("A17").Value = "zine"

# The "tx-sucesso" sheet also needs the two tied rows swapped.
$wsTx = $wb.Worksheets.Item("tx-sucesso")

$wsTx.Range("A8").Value = "questoes_genero"
$wsTx.Range("A9").Value = "erotismo"
